# close #206: Adds support for zero-sum influencing factors
#
# Two "influencing factor" pairs (P4/Q4 and P5/Q5) are reworked so a
# factor that nets out to zero is shown as the literal marker "DI"
# ("Diferencia indeterminada / N/A") instead of a (meaningless) computed
# ratio, while its paired column collapses to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (id 1100049): P4 becomes the zero-sum value, Q4 becomes the "DI" marker.
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = "DI"

# Row 5 (id 1100056): mirror image - P5 becomes "DI", Q5 becomes the zero-sum value.
$ws.Range("P5").Value = "DI"
$ws.Range("Q5").Value = 0

# Reflect where the analyst ended up looking after making the change:
# scrolled the frozen-header view right to column H, then left the
# selection on Q22 (just past the data, bottom-right corner area).
$ws.Range("H1").Select() | Out-Null
$ws.Range("Q22").Select() | Out-Null
